# This workbook tracks weekly Ají (chili pepper) price records.
# The edit inserts 3 new price records (one weekly update) right before the
# existing row 212, pushing the rows that used to be 212-308 down to 215-311.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows before row 212; everything from 212..308 shifts to 215..311.
$ws.Rows("212:214").Insert()

# Columns that stay constant across every data row of this sheet
# (A=Mercado ID, B=Mercado, C=Region, E=Codreg, F=Categoria ID, G=Categoria,
#  N=Unidad de comercializacion, O=Origen, Q=Kg o Unidades, R=Clasificacion).
# Copy them from the row right below (215, the former row 212) into the
# newly inserted rows so the new records carry the same fixed metadata.
$commonCols = @(1, 2, 3, 5, 6, 7, 14, 15, 17, 18)
for ($r = 212; $r -le 214; $r++) {
    foreach ($col in $commonCols) {
        $ws.Cells.Item($r, $col).Value = $ws.Cells.Item(215, $col).Value2
    }
}

# Row 212: new Ají "Americana (o)" / "Primera" record
$ws.Cells.Item(212, 4).Value = 44784   # D Fecha
$ws.Cells.Item(212, 8).Value = "Americana (o)"   # H Variedad
$ws.Cells.Item(212, 9).Value = "Primera"         # I Calidad
$ws.Cells.Item(212, 10).Value = 60     # J Volumen
$ws.Cells.Item(212, 11).Value = 40000  # K Precio minimo
$ws.Cells.Item(212, 12).Value = 42000  # L Precio maximo
$ws.Cells.Item(212, 13).Value = 41000  # M Precio promedio ponderado
$ws.Cells.Item(212, 16).Value = 1640   # P Precio $/Kg

# Row 213: new Ají "Americana (o)" / "Segunda" record
$ws.Cells.Item(213, 4).Value = 44784
$ws.Cells.Item(213, 8).Value = "Americana (o)"
$ws.Cells.Item(213, 9).Value = "Segunda"
$ws.Cells.Item(213, 10).Value = 40
$ws.Cells.Item(213, 11).Value = 30000
$ws.Cells.Item(213, 12).Value = 32000
$ws.Cells.Item(213, 13).Value = 31000
$ws.Cells.Item(213, 16).Value = 1240

# Row 214: new Ají "Inferno" / "Primera" record
$ws.Cells.Item(214, 4).Value = 44784
$ws.Cells.Item(214, 8).Value = "Inferno"
$ws.Cells.Item(214, 9).Value = "Primera"
$ws.Cells.Item(214, 10).Value = 40
$ws.Cells.Item(214, 11).Value = 25000
$ws.Cells.Item(214, 12).Value = 26000
$ws.Cells.Item(214, 13).Value = 25500
$ws.Cells.Item(214, 16).Value = 1020
